# Change AddressBook to HealthBook for command sequence diagrams.
#
# 1) On the (single) content slide, two shapes reference "AddressBook":
#      - the "AddressBookParser" participant box -> "HealthParser"... no,
#        actually only the "Address" part of ":AddressBookParser" changes
#        to ":Health" (the second paragraph "BookParser" is untouched).
#      - the ":VersionedAddressBook" participant box -> ":VersionedHealthBook"
# 2) Every "24 Oct 2018" auto-date placeholder (all 11 slide layouts, the
#    slide master, and the notes master) is bumped to "12 Nov 2018".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide text updates
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($shp.Id -eq 9) {
        # ":AddressBookParser" (paragraph 1 = ":Address", paragraph 2 = "BookParser")
        # Rewrite paragraph 1 in one go so it stays a single run, like the source.
        $para1 = $tr.Characters(1, 8)
        if ($para1.Text -eq ":Address") {
            $para1.Text = ":Health"
        }
    }
    elseif ($shp.Id -eq 26) {
        # ":VersionedAddressBook" -> ":VersionedHealthBook"
        # Leave the leading ":" run untouched; only rewrite the second run.
        $target = $tr.Characters(2, 20)
        if ($target.Text -eq "VersionedAddressBook") {
            $target.Text = "VersionedHealthBook"
        }
    }
}

# ---------------------------------------------------------------------
# 2. Date placeholder updates ("24 Oct 2018" -> "12 Nov 2018")
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $s = $shapes.Item($k)
        if (-not $s.HasTextFrame) { continue }
        $t = $s.TextFrame.TextRange
        if ($t.Text -eq "24 Oct 2018") {
            $t.Text = "12 Nov 2018"
        }
    }
}

$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

Update-DateShapes $master.Shapes

$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
